# Apply "Batterywise analysis" relabeling / value update edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 6-10: relabel + swap/update values ---
$ws.Range("A6").Value = "Starting SoC (%)"
$ws.Range("B6").Value = 100

$ws.Range("A7").Value = "Ending SoC (%)"
$ws.Range("B7").Value = 15

$ws.Range("A8").Value = "Total distance covered (km)"

$ws.Range("A9").Value = "Total energy consumption(WH/KM)"

$ws.Range("A10").Value = "Total SOC consumed(%)"
$ws.Range("B10").Value = 85

# --- Rows 12-30: append unit suffixes, a couple value updates, and a swap ---
$ws.Range("A12").Value = "Peak Power(kW)"
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"

$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 4.047046564394381

$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.483

$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.144

$ws.Range("A18").Value = "Difference in Cell Voltage(V)"
$ws.Range("A19").Value = "Minimum Temperature(C)"
$ws.Range("A20").Value = "Maximum Temperature(C)"

$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 12

$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("B28").Value = 48

$ws.Range("A29").Value = "lowest cell temp(C)"
$ws.Range("B29").Value = 36

$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# --- Row 31 ("Maximum BMS Temperature in C") is removed; rows 32-42 shift up to 31-41,
#     relabeled with unit suffixes and new values; a new row 43 is appended. ---
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 56

$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.743151371111111

$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.00000007264946949700389

$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 7.10488322717622

$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 5.470912951167728

$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 8.798301486199575

$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 14.21656050955414

$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 41.97537154989384

$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 21.1176220806794

$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 0

$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 0

$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 0

$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
